$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 250.57143
$ws.Range("I2").Value = 241.9
$ws.Range("K2").Value = 241.9
$ws.Range("M2").Value = -128.9
# Row 33
$ws.Range("H33").Value = 193.42105
$ws.Range("I33").Value = 155.7
$ws.Range("J33").Value = 235.33333
$ws.Range("K33").Value = 155.7
$ws.Range("L33").Value = 235.33333
$ws.Range("M33").Value = 73.30000000000001
$ws.Range("N33").Value = -693.3333299999999
# Row 34
$ws.Range("H34").Value = 5604
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()
# Row 36
$ws.Range("H36").Value = 5604
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()
# Row 46
$ws.Range("H46").Value = 4998
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()
# Row 55
$ws.Range("H55").Value = 536.8
$ws.Range("I55").Value = 73.8
$ws.Range("J55").Value = 999.8
$ws.Range("K55").Value = 73.8
$ws.Range("L55").Value = 999.8
$ws.Range("M55").Value = 140.2
$ws.Range("N55").Value = -1427.8
# Row 60
$ws.Range("H60").Value = 4998
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()
# Row 101
$ws.Range("H101").Value = 2165
$ws.Range("I101").Value = 2431
$ws.Range("K101").Value = 7293
$ws.Range("M101").Value = -5671

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 37
$ws.Range("H37").Value = 24000
$ws.Range("I37").Value = 24000
$ws.Range("K37").Value = 24000
$ws.Range("M37").Value = -23727
# Row 61
$ws.Range("H61").Value = 4481.1665
$ws.Range("I61").Value = 3972.25
$ws.Range("K61").Value = 3972.25
$ws.Range("M61").Value = -3760.25
# Row 63
$ws.Range("H63").Value = 4573.1
$ws.Range("I63").Value = 1373.75
$ws.Range("J63").Value = 6706
$ws.Range("K63").Value = 1373.75
$ws.Range("L63").Value = 6706
$ws.Range("M63").Value = -687.75
$ws.Range("N63").Value = -8078
# Row 66
$ws.Range("H66").Value = 4573.1
$ws.Range("I66").Value = 1373.75
$ws.Range("J66").Value = 6706
$ws.Range("K66").Value = 6868.75
$ws.Range("L66").Value = 33530
$ws.Range("M66").Value = -3436.75
$ws.Range("N66").Value = -40394
# Row 74
$ws.Range("H74").Value = 689.2222
$ws.Range("I74").Value = 689.2222
$ws.Range("K74").Value = 689.2222
$ws.Range("M74").Value = 184.7778
# Row 77
$ws.Range("H77").Value = 689.2222
$ws.Range("I77").Value = 689.2222
$ws.Range("K77").Value = 3446.111
$ws.Range("M77").Value = 921.8889999999997
# Row 136
$ws.Range("H136").Value = 4481.1665
$ws.Range("I136").Value = 3972.25
$ws.Range("K136").Value = 11916.75
$ws.Range("M136").Value = -9366.75

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 2147
$ws.Range("I20").Value = 1908.625
$ws.Range("K20").Value = 1908.625
$ws.Range("M20").Value = -1661.625
# Row 29
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()
# Row 36
$ws.Range("H36").Value = 11333
$ws.Range("I36").Value = 3000
$ws.Range("J36").Value = 27999
$ws.Range("K36").Value = 3000
$ws.Range("L36").Value = 27999
$ws.Range("M36").Value = -2466
$ws.Range("N36").Value = -29067
# Row 95
$ws.Range("H95").Value = 22651.75
$ws.Range("J95").Value = 22651.75
$ws.Range("L95").Value = 22651.75
$ws.Range("N95").Value = -28143.75
# Row 134
$ws.Range("H134").Value = 6255.5
$ws.Range("I134").Value = 6157.7
$ws.Range("K134").Value = 18473.1
$ws.Range("M134").Value = -15938.1

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 1346.8
$ws.Range("I22").Value = 694.5
$ws.Range("J22").Value = 1781.6666
$ws.Range("K22").Value = 694.5
$ws.Range("L22").Value = 1781.6666
$ws.Range("M22").Value = -344.5
$ws.Range("N22").Value = -2481.6666
# Row 31
$ws.Range("H31").Value = 2708.2666
$ws.Range("I31").Value = 2086.5386
$ws.Range("K31").Value = 2086.5386
$ws.Range("M31").Value = -1791.5386
# Row 34
$ws.Range("H34").Value = 2708.2666
$ws.Range("I34").Value = 2086.5386
$ws.Range("K34").Value = 2086.5386
$ws.Range("M34").Value = -1884.5386
# Row 55
$ws.Range("H55").Value = 8000
$ws.Range("I55").Value = 8000
$ws.Range("K55").Value = 8000
$ws.Range("M55").Value = -7685
# Row 58
$ws.Range("H58").Value = 3996
$ws.Range("I58").Value = 1839.8
$ws.Range("J58").Value = 7589.6665
$ws.Range("K58").Value = 1839.8
$ws.Range("L58").Value = 7589.6665
$ws.Range("M58").Value = -1636.8
$ws.Range("N58").Value = -7995.6665
# Row 133
$ws.Range("H133").Value = 80000
$ws.Range("J133").Value = 80000
$ws.Range("L133").Value = 80000
$ws.Range("N133").Value = -85060
# Row 136
$ws.Range("H136").Value = 3996
$ws.Range("I136").Value = 1839.8
$ws.Range("J136").Value = 7589.6665
$ws.Range("K136").Value = 5519.4
$ws.Range("L136").Value = 22768.9995
$ws.Range("M136").Value = -2969.4
$ws.Range("N136").Value = -27868.9995

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 98
$ws.Range("H98").Value = 717
$ws.Range("J98").Value = 717
$ws.Range("L98").Value = 2151
$ws.Range("N98").Value = -5147

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 111
$ws.Range("H111").Value = 29500
$ws.Range("J111").Value = 29500
$ws.Range("L111").Value = 29500
$ws.Range("N111").Value = -35634
# Row 132
$ws.Range("H132").Value = 2578.3044
$ws.Range("I132").Value = 2348.7144
$ws.Range("J132").Value = 4989
$ws.Range("K132").Value = 7046.1432
$ws.Range("L132").Value = 14967
$ws.Range("M132").Value = -4516.1432
$ws.Range("N132").Value = -20027

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 300
$ws.Range("I16").Value = 300
$ws.Range("K16").Value = 300
$ws.Range("M16").Value = -130
# Row 29
$ws.Range("H29").Value = 24999.5
$ws.Range("I29").Value = 24999.5
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 24999.5
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -24704.5
$ws.Range("N29").ClearContents()
# Row 31
$ws.Range("H31").Value = 4105.7144
$ws.Range("I31").Value = 4105.7144
$ws.Range("K31").Value = 4105.7144
$ws.Range("M31").Value = -3857.7144
# Row 35
$ws.Range("H35").Value = 665.625
$ws.Range("I35").Value = 665.625
$ws.Range("K35").Value = 665.625
$ws.Range("M35").Value = -329.625
# Row 68
$ws.Range("H68").Value = 1375
# Row 71
$ws.Range("H71").Value = 1375
# Row 122
$ws.Range("H122").Value = 6247.1665
$ws.Range("I122").Value = 5997.1
$ws.Range("J122").Value = 7497.5
$ws.Range("K122").Value = 17991.3
$ws.Range("L122").Value = 22492.5
$ws.Range("M122").Value = -15541.3
$ws.Range("N122").Value = -27392.5
# Row 136
$ws.Range("H136").Value = 3975.3845
$ws.Range("J136").Value = 4999.5
$ws.Range("L136").Value = 14998.5
$ws.Range("N136").Value = -20098.5

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 418.22223
$ws.Range("J107").Value = 333.66666
$ws.Range("L107").Value = 1000.99998
$ws.Range("N107").Value = -4840.99998
# Row 132
$ws.Range("H132").Value = 2237.7
$ws.Range("I132").Value = 1922.125
$ws.Range("K132").Value = 5766.375
$ws.Range("M132").Value = -3236.375
